$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A plain, unstyled data cell used as a formatting donor (General number format,
# default style) so that cells which need an explicit Text format to avoid
# numeric auto-conversion can have their style reset back to the workbook default
# after the value is written.
$styleDonor = $ws.Range("A2")

$ws.Range("D2").Value = '57.289.37'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '3.094.21'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.91'
$styleDonor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.54'
$styleDonor.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$styleDonor.Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '3.093.38'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.458'
$styleDonor.Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.36'
$styleDonor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = '  +2.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.107'
$styleDonor.Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.401'
$styleDonor.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '3.625.26'
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.136'
$styleDonor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.39'
$styleDonor.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000162'
$styleDonor.Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").Value = '57.391.97'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '3.090.50'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.90'
$styleDonor.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = '  -3.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.49'
$styleDonor.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.88'
$styleDonor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = '  -1.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '350.41'
$styleDonor.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.78'
$styleDonor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = '  +1.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.499'
$styleDonor.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$styleDonor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = '  -1.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$styleDonor.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("D28").Value = '0.0₃0866'
$ws.Range("E28").Value = '  -6.40%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.87'
$styleDonor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.87'
$styleDonor.Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = '  -7.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.94'
$styleDonor.Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.87'
$styleDonor.Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = '  +5.30%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.14'
$styleDonor.Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("E35").Value = '  -3.25%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.17'
$styleDonor.Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.03'
$styleDonor.Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.64'
$styleDonor.Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0659'
$styleDonor.Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.06'
$styleDonor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.695'
$styleDonor.Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").Value = '2.393.55'
$ws.Range("E44").Value = '  +5.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.72'
$styleDonor.Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '3.134.33'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0264'
$styleDonor.Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.958'
$styleDonor.Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = '  -3.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.96'
$styleDonor.Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.762'
$styleDonor.Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = '  +1.48%  '
